$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-06-23 Sunday" "2024-06-24 Monday"
Replace-Text "254×5=" "617×3="
Replace-Text "628×9=" "833×8="
Replace-Text "806×9=" "166×2="
Replace-Text "991×9=" "407×3="
Replace-Text "200×3=" "853×8="
Replace-Text "942×2=" "863×3="
Replace-Text "704×3=" "138×8="
Replace-Text "347×9=" "499×3="
Replace-Text "438×4=" "854×9="
Replace-Text "422×4=" "207×6="
Replace-Text "114×4=" "490×5="
Replace-Text "758×2=" "703×4="
Replace-Text "518×2=" "607×8="
Replace-Text "194×9=" "510×4="
Replace-Text "828×8=" "777×8="
Replace-Text "260×5=" "522×5="
Replace-Text "135×4=" "286×8="
Replace-Text "627×9=" "946×7="
Replace-Text "644×4=" "116×3="
Replace-Text "251×6=" "155×9="
Replace-Text "820×7=" "153×3="
Replace-Text "609×9=" "330×8="
Replace-Text "674×7=" "311×4="
Replace-Text "293×7=" "770×4="
Replace-Text "163×8=" "999×3="
